$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header of column B first (introduces new shared string before others)
$ws.Range("B1").Value = "Islands without Magnis and Coriorrostris"

# Column B: replace "Espanola" (row4) with shifted-up values, append two new entries
$ws.Range("B4").Value = "San Cristobal"
$ws.Range("B5").Value = "Santa Fe"
$ws.Range("B6").Value = "Champion"
$ws.Range("B7").Value = "Baltra"
$ws.Range("B8").Value = "Enderby"
$ws.Range("B9").Value = "Gardner"
$ws.Range("B10").Value = "Daphne Major <1983"

# Column A: add "Espanola" as new row 16
$ws.Range("A16").Value = "Espanola"

# Update header of column A last
$ws.Range("A1").Value = "Islands with Magnis and Cornirrostris"

# Update dimension-driving selection (cosmetic, mirrors typical post-edit state)
$ws.Range("E28").Select()

# Re-fit column widths to match new (longer) content, as Excel does automatically
$ws.Columns.Item(1).ColumnWidth = 33.5
$ws.Columns.Item(2).ColumnWidth = 36.5
